$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.347.28'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.846.22'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '239.92'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").Value = '0.6269'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").Value = '0.9985'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '0.07606'
$ws.Range("E8").Value = '  -0.83%  '
$ws.Range("D9").Value = '0.2899'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("D10").Value = '24.70'
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").Value = '5.023'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '0.6781'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '0.00001050'
$ws.Range("E14").Value = '  -3.67%  '
$ws.Range("D15").Value = '82.94'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '6.145'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '29.368.13'
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("D18").Value = '227.38'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").Value = '12.32'
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("D20").Value = '0.9984'
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").Value = '7.463'
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").Value = '0.9984'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '158.18'
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").Value = '0.1384'
$ws.Range("E24").Value = '  -0.42%  '
$ws.Range("D25").Value = '8.398'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = '17.65'
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").Value = '1.404'
$ws.Range("E27").Value = '  +7.03%  '
$ws.Range("D28").Value = '1.460'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").Value = '0.05594'
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").Value = '4.102'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").Value = '4.058'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").Value = '1.162'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = '1.833'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '0.6967'
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").Value = '2.585'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").Value = '0.01799'
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.227.63'
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").Value = '2.724'
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("D39").Value = '6.369'
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("D40").Value = '0.9007'
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("D41").Value = '0.9986'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '101.18'
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("D43").Value = '65.67'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '7.206'
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000119'
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.3991'
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '8.991'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.675'
$ws.Range("E48").Value = '  -0.83%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1139'
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.05701'
$ws.Range("E50").Value = '  -0.22%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.4620'
$ws.Range("E51").Value = '  -0.18%  '
